$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price/volume table (Price in column D, Volume(1h) in
# column E) with newly scraped figures. Price cells that look like plain
# numbers (e.g. "36.68") are written with a leading apostrophe so Excel
# keeps them as text instead of silently converting to a Number, then the
# style is reset to "Normal" so the quote-prefix indicator introduced by
# the apostrophe doesn't leave a stray format behind.

$ws.Range("D2").Value = "63.513.99"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "3.070.07"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'591.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'154.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "3.071.61"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").Value = "'5.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").Value = "'0.0000237"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").Value = "'36.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "3.582.34"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "'7.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "63.496.10"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").Value = "3.072.68"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "'484.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.02%  "
$ws.Range("D21").Value = "'14.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "'0.709"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.16%  "
$ws.Range("D23").Value = "'7.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "'2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").Value = "'81.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "'12.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.96%  "
$ws.Range("D27").Value = "'10.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.02%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "'7.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").Value = "'2.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("D34").Value = "'27.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "0.0₃0824"
$ws.Range("E36").Value = "  -3.19%  "
$ws.Range("D37").Value = "'6.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").Value = "'3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.42%  "
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").Value = "'50.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'443.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("D46").Value = "2.824.64"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").Value = "'132.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'25.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  -2.16%  "
